$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics for Cd47-Sirpa LR pairs
# Row 2
$ws.Range("G2").Value = 33.54649666666666
$ws.Range("H2").Value = 100.63949
$ws.Range("I2").Value = 0.1978943147725085
$ws.Range("J2").Value = 0.1978943147725085
$ws.Range("M2").Value = 0.197995
$ws.Range("N2").Value = 0.593985
$ws.Range("O2").Value = 0.00464724904661508
$ws.Range("P2").Value = 0.00464724904661508
$ws.Range("Q2").Value = 6.642038607516666
$ws.Range("R2").Value = 59.77834746764999
$ws.Range("S2").Value = 0.0009196641656570847
$ws.Range("T2").Value = 0.0009196641656570847

# Row 3
$ws.Range("G3").Value = 33.54649666666666
$ws.Range("H3").Value = 100.63949
$ws.Range("I3").Value = 0.1978943147725085
$ws.Range("J3").Value = 0.1978943147725085
$ws.Range("O3").Value = 0.8589674670659959
$ws.Range("P3").Value = 0.8589674670659959
$ws.Range("Q3").Value = 1227.671472224781
$ws.Range("R3").Value = 11049.04325002303
$ws.Range("S3").Value = 0.1699847783069025
$ws.Range("T3").Value = 0.1699847783069025

# Row 4
$ws.Range("G4").Value = 33.54649666666666
$ws.Range("H4").Value = 100.63949
$ws.Range("I4").Value = 0.1978943147725085
$ws.Range("J4").Value = 0.1978943147725085
$ws.Range("O4").Value = 0.136385283887389
$ws.Range("P4").Value = 0.136385283887389
$ws.Range("Q4").Value = 194.9274316892855
$ws.Range("R4").Value = 1754.34688520357
$ws.Range("S4").Value = 0.02698987229994889
$ws.Range("T4").Value = 0.02698987229994889

# Row 5
$ws.Range("I5").Value = 0.6322484766686425
$ws.Range("J5").Value = 0.6322484766686425
$ws.Range("M5").Value = 0.197995
$ws.Range("N5").Value = 0.593985
$ws.Range("O5").Value = 0.00464724904661508
$ws.Range("P5").Value = 0.00464724904661508
$ws.Range("Q5").Value = 21.22051255693833
$ws.Range("R5").Value = 190.984613012445
$ws.Range("S5").Value = 0.002938216130422185
$ws.Range("T5").Value = 0.002938216130422185

# Row 6
$ws.Range("I6").Value = 0.6322484766686425
$ws.Range("J6").Value = 0.6322484766686425
$ws.Range("O6").Value = 0.8589674670659959
$ws.Range("P6").Value = 0.8589674670659959
$ws.Range("S6").Value = 0.5430808725603983
$ws.Range("T6").Value = 0.5430808725603983

# Row 7
$ws.Range("I7").Value = 0.6322484766686425
$ws.Range("J7").Value = 0.6322484766686425
$ws.Range("O7").Value = 0.136385283887389
$ws.Range("P7").Value = 0.136385283887389
$ws.Range("S7").Value = 0.08622938797782202
$ws.Range("T7").Value = 0.08622938797782204

# Row 8
$ws.Range("I8").Value = 0.169857208558849
$ws.Range("J8").Value = 0.169857208558849
$ws.Range("M8").Value = 0.197995
$ws.Range("N8").Value = 0.593985
$ws.Range("O8").Value = 0.00464724904661508
$ws.Range("P8").Value = 0.00464724904661508
$ws.Range("Q8").Value = 5.701013383379999
$ws.Range("R8").Value = 51.30912045042
$ws.Range("S8").Value = 0.0007893687505358099
$ws.Range("T8").Value = 0.0007893687505358096

# Row 9
$ws.Range("I9").Value = 0.169857208558849
$ws.Range("J9").Value = 0.169857208558849
$ws.Range("O9").Value = 0.8589674670659959
$ws.Range("P9").Value = 0.8589674670659959
$ws.Range("S9").Value = 0.1459018161986951
$ws.Range("T9").Value = 0.1459018161986951

# Row 10
$ws.Range("I10").Value = 0.169857208558849
$ws.Range("J10").Value = 0.169857208558849
$ws.Range("O10").Value = 0.136385283887389
$ws.Range("P10").Value = 0.136385283887389
$ws.Range("S10").Value = 0.02316602360961805
$ws.Range("T10").Value = 0.02316602360961806
